# Regen save_data to use K (strikeouts) instead of Strike# in column G.
# Recalculated K values are written back into G2:G73 (the "K" column),
# leaving every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2..73 (header "K" lives in G1).
$kValues = @(2,1,1,2,1,2,0,1,0,1,1,0,1,1,1,2,3,1,1,0,1,2,2,0,2,1,1,1,3,1,2,3,2,1,2,1,1,1,0,1,1,0,2,1,1,2,0,1,3,2,0,0,2,0,2,0,1,1,3,0,2,1,1,0,1,1,4,1,0,1,2,0)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
